$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3595.7646
$ws.Range("J17").Value = 2933.4666
$ws.Range("L17").Value = 8800.399800000001
$ws.Range("N17").Value = -9136.399800000001

$ws.Range("H96").Value = 1296.2778
$ws.Range("I96").Value = 1042.625
$ws.Range("J96").Value = 1499.2
$ws.Range("K96").Value = 3127.875
$ws.Range("L96").Value = 4497.6
$ws.Range("M96").Value = -1754.875
$ws.Range("N96").Value = -7243.6

$ws.Range("H98").Value = 1469.6
$ws.Range("I98").Value = 1073.9565
$ws.Range("J98").Value = 2769.5715
$ws.Range("K98").Value = 1073.9565
$ws.Range("L98").Value = 2769.5715
$ws.Range("M98").Value = 424.0435
$ws.Range("N98").Value = -5765.5715

$ws.Range("H112").Value = 1458.6923
$ws.Range("J112").Value = 1458.6923
$ws.Range("L112").Value = 4376.0769
$ws.Range("N112").Value = -6592.0769

$ws.Range("H122").Value = 1469.6
$ws.Range("I122").Value = 1073.9565
$ws.Range("J122").Value = 2769.5715
$ws.Range("K122").Value = 3221.8695
$ws.Range("L122").Value = 8308.7145
$ws.Range("M122").Value = -771.8694999999998
$ws.Range("N122").Value = -13208.7145

$ws.Range("H132").Value = 857.2162
$ws.Range("I132").Value = 745.89655
$ws.Range("K132").Value = 2237.68965
$ws.Range("M132").Value = 292.3103499999997

$ws.Range("H137").Value = 2749.1428
$ws.Range("I137").Value = 2665.4443
$ws.Range("K137").Value = 7996.3329
$ws.Range("M137").Value = -5446.3329

$ws.Range("H138").Value = 1911.254
$ws.Range("I138").Value = 1551.7097
$ws.Range("J138").Value = 2259.5625
$ws.Range("K138").Value = 4655.1291
$ws.Range("L138").Value = 6778.6875
$ws.Range("M138").Value = 484.8708999999999
$ws.Range("N138").Value = -17058.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1205.7037
$ws.Range("J74").Value = 2650.5557
$ws.Range("L74").Value = 2650.5557
$ws.Range("N74").Value = -4398.5557

$ws.Range("H77").Value = 1205.7037
$ws.Range("J77").Value = 2650.5557
$ws.Range("L77").Value = 13252.7785
$ws.Range("N77").Value = -21988.7785

$ws.Range("H104").Value = 38800
$ws.Range("J104").Value = 39750
$ws.Range("L104").Value = 39750
$ws.Range("N104").Value = -46738

$ws.Range("H109").Value = 61748
$ws.Range("J109").Value = 61748
$ws.Range("L109").Value = 61748
$ws.Range("N109").Value = -64522

$ws.Range("H123").Value = 82000
$ws.Range("J123").Value = 82000
$ws.Range("L123").Value = 82000
$ws.Range("N123").Value = -91800

$ws.Range("H132").Value = 2188.6843
$ws.Range("I132").Value = 1779.5454
$ws.Range("J132").Value = 2751.25
$ws.Range("K132").Value = 5338.6362
$ws.Range("L132").Value = 8253.75
$ws.Range("M132").Value = -2808.6362
$ws.Range("N132").Value = -13313.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5399.75
$ws.Range("I134").Value = 6393.8096
$ws.Range("K134").Value = 19181.4288
$ws.Range("M134").Value = -16646.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 153.33333
$ws.Range("I7").Value = 153.33333
$ws.Range("K7").Value = 153.33333
$ws.Range("M7").Value = -40.33332999999999

$ws.Range("H9").Value = 17500
$ws.Range("J9").Value = 17500
$ws.Range("L9").Value = 17500
$ws.Range("N9").Value = -17836

$ws.Range("H31").Value = 2104.5454
$ws.Range("I31").Value = 1733.3334
$ws.Range("J31").Value = 2243.75
$ws.Range("K31").Value = 1733.3334
$ws.Range("L31").Value = 2243.75
$ws.Range("M31").Value = -1438.3334
$ws.Range("N31").Value = -2833.75

$ws.Range("H34").Value = 2104.5454
$ws.Range("I34").Value = 1733.3334
$ws.Range("J34").Value = 2243.75
$ws.Range("K34").Value = 1733.3334
$ws.Range("L34").Value = 2243.75
$ws.Range("M34").Value = -1531.3334
$ws.Range("N34").Value = -2647.75

$ws.Range("H58").Value = 3346387.2
$ws.Range("I58").Value = 6211771.5
$ws.Range("J58").Value = 3438.6667
$ws.Range("K58").Value = 6211771.5
$ws.Range("L58").Value = 3438.6667
$ws.Range("M58").Value = -6211568.5
$ws.Range("N58").Value = -3844.6667

$ws.Range("H136").Value = 3346387.2
$ws.Range("I136").Value = 6211771.5
$ws.Range("J136").Value = 3438.6667
$ws.Range("K136").Value = 18635314.5
$ws.Range("L136").Value = 10316.0001
$ws.Range("M136").Value = -18632764.5
$ws.Range("N136").Value = -15416.0001

$ws.Range("H141").Value = 63798.8
$ws.Range("J141").Value = 61748.5
$ws.Range("L141").Value = 61748.5
$ws.Range("N141").Value = -72108.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 285.875
$ws.Range("I2").Value = 272.5
$ws.Range("J2").Value = 299.25
$ws.Range("K2").Value = 1635
$ws.Range("L2").Value = 1795.5
$ws.Range("M2").Value = -1522
$ws.Range("N2").Value = -2021.5

$ws.Range("H56").Value = 6448.476
$ws.Range("I56").Value = 6448.476
$ws.Range("K56").Value = 6448.476
$ws.Range("M56").Value = -5918.476

$ws.Range("H75").Value = 26316.2
$ws.Range("J75").Value = 26316.2
$ws.Range("L75").Value = 78948.60000000001
$ws.Range("N75").Value = -80944.60000000001

$ws.Range("H78").Value = 26316.2
$ws.Range("J78").Value = 26316.2
$ws.Range("L78").Value = 236845.8
$ws.Range("N78").Value = -246829.8

$ws.Range("H117").Value = 609.9
$ws.Range("I117").Value = 425.75
$ws.Range("J117").Value = 732.6667
$ws.Range("K117").Value = 1277.25
$ws.Range("L117").Value = 2198.0001
$ws.Range("M117").Value = 2164.75
$ws.Range("N117").Value = -9082.000100000001

$ws.Range("H131").Value = 11564.2295
$ws.Range("J131").Value = 12026.662
$ws.Range("L131").Value = 36079.986
$ws.Range("N131").Value = -46159.986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6728.4
$ws.Range("I16").Value = 8223
$ws.Range("K16").Value = 8223
$ws.Range("M16").Value = -8053

$ws.Range("H46").Value = 1887.6923
$ws.Range("I46").Value = 1095.125
$ws.Range("K46").Value = 1095.125
$ws.Range("M46").Value = -907.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 14666.667
$ws.Range("J21").Value = 14666.667
$ws.Range("L21").Value = 14666.667
$ws.Range("N21").Value = -15136.667

$ws.Range("H35").Value = 14666.667
$ws.Range("J35").Value = 14666.667
$ws.Range("L35").Value = 14666.667
$ws.Range("N35").Value = -15246.667

$ws.Range("H100").Value = 1197.2
$ws.Range("I100").Value = 1024.5714
$ws.Range("K100").Value = 2049.1428
$ws.Range("M100").Value = -1508.1428

$ws.Range("H122").Value = 33229.168
$ws.Range("I122").Value = 37676.285
$ws.Range("K122").Value = 113028.855
$ws.Range("M122").Value = -110578.855

$ws.Range("H123").Value = 48055.727
$ws.Range("J123").Value = 48055.727
$ws.Range("L123").Value = 48055.727
$ws.Range("N123").Value = -57855.727

$ws.Range("H132").Value = 1132.4603
$ws.Range("I132").Value = 890.087
$ws.Range("J132").Value = 1788.2941
$ws.Range("K132").Value = 2670.261
$ws.Range("L132").Value = 5364.8823
$ws.Range("M132").Value = -140.261
$ws.Range("N132").Value = -10424.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1551
$ws.Range("I136").Value = 1551
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4653
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 447
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 12000
$ws.Range("J24").Value = 12000
$ws.Range("L24").Value = 12000
$ws.Range("N24").Value = -12460

$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
